# repull data, push all data, mean calculation
# Update column F ("dSF") values for the rows whose pulled data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    3  = 3
    6  = 2
    7  = 2
    8  = 3
    9  = -2
    11 = 7
    12 = 1
    13 = -2
    14 = 5
    16 = 2
    17 = -2
    18 = 7
    19 = 5
    20 = 2
    21 = 3
    22 = -2
    23 = 3
    24 = 3
    25 = 4
    26 = -2
    27 = 1
    28 = 1
    29 = 2
    30 = 2
    31 = -1
    32 = -6
    33 = 3
    34 = 3
    36 = -2
    37 = 2
    38 = -3
    40 = 1
    41 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
